$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three idiom/story rows that were removed from the workbook:
#   - row 382 : 模棱两可
#   - row 684 : 一笔勾销  (becomes row 683 after the first deletion shifts rows up)
#   - row 764 : 予取予求  (becomes row 762 after the first two deletions shift rows up)
# Deleting from top to bottom and adjusting indices for the upward shift caused
# by each prior deletion reproduces the same net result as deleting all three
# in a single pass from the bottom up.
$ws.Rows(382).Delete()
$ws.Rows(683).Delete()
$ws.Rows(762).Delete()
